$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 768, shifting existing rows 768-846 down to 769-847
$ws.Rows(768).Insert()

# Populate the newly inserted row 768 with the new data entry
$ws.Range("A768").Value = 3
$ws.Range("B768").Value = "Femacal de La Calera"
$ws.Range("C768").Value = "Coquimbo"
$ws.Range("D768").Value = 45194
$ws.Range("E768").Value = 5
$ws.Range("F768").Value = 100112032
$ws.Range("G768").Value = "Zapallo italiano"
$ws.Range("H768").Value = "Sin especificar"
$ws.Range("I768").Value = "Primera"
$ws.Range("J768").Value = 80
$ws.Range("K768").Value = 13000
$ws.Range("L768").Value = 14000
$ws.Range("M768").Value = 13500
$ws.Range("N768").Value = "`$/caja 60 unidades"
$ws.Range("O768").Value = "Región de Arica y Parinacota"
$ws.Range("P768").Value = 225
$ws.Range("Q768").Value = 60
$ws.Range("R768").Value = "Hortaliza"
